# Apply cyclic update to rows 2-8 on the "Artfynd" sheet.
# Each row's Id (A), Ost (Q), Nord (R), Startdatum (Y) and Slutdatum (AA)
# values are replaced by the values that originally belonged to the next
# row down (row 8 receives the values that originally belonged to row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# New values to write into rows 2-8 (columns A, Q, R, Y, AA)
$newValues = @(
    @{ Row = 2; A = 94726850; Q = 385883.915955339;  R = 6857995.776764749; Date = "2021-06-08" }
    @{ Row = 3; A = 94726809; Q = 386572.4667873504; R = 6857676.512812085; Date = "2021-06-15" }
    @{ Row = 4; A = 94726815; Q = 386513.1098144559; R = 6857970.586442914; Date = "2021-06-11" }
    @{ Row = 5; A = 94726808; Q = 386471.5539895515; R = 6857773.923141795; Date = "2021-06-15" }
    @{ Row = 6; A = 94726849; Q = 386436.8180425634; R = 6858381.457110016; Date = "2021-06-09" }
    @{ Row = 7; A = 94726810; Q = 386473.1449506464; R = 6857779.069685961; Date = "2021-06-15" }
    @{ Row = 8; A = 94726847; Q = 386462.6085887029; R = 6858245.91119291;  Date = "2021-06-10" }
)

# The Startdatum/Slutdatum columns hold plain-text dates (e.g. "2021-06-08")
# in the original workbook. Assigning a date-shaped string directly to
# Range.Value would make Excel auto-convert it into a date serial number,
# so force a text number format on those columns first, then restore the
# default ("Normal") style once the text has been written so no stray
# formatting is left behind on the cells. (Using two single-column ranges
# instead of one multi-area range avoids the second area silently keeping
# its default/general number format.)
$yRange = $ws.Range("Y2:Y8")
$aaRange = $ws.Range("AA2:AA8")
$yRange.NumberFormat = "@"
$aaRange.NumberFormat = "@"

foreach ($entry in $newValues) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.A
    $ws.Range("Q$r").Value = $entry.Q
    $ws.Range("R$r").Value = $entry.R
    $ws.Range("Y$r").Value = $entry.Date
    $ws.Range("AA$r").Value = $entry.Date
}

$yRange.Style = "Normal"
$aaRange.Style = "Normal"
